# Apply "remove duplicates" edit to df_mini sheet.
# Net effect (verified against the canonical XML diff):
#   - Column D values for rows 13-16 and 18-21 shift up by one (duplicate
#     "Clubman" rows removed from the original data), while columns A, C, E
#     stay exactly as they were for each row position.
#   - Row 17 additionally gets its Year (B) and Additional Info (F) updated
#     to match the row that "moved into" position 17.
#   - The two trailing rows (22 and 23), which become redundant once the
#     duplicates are squeezed out, are deleted entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column D (Model) for rows 13-16 ---
$ws.Range("D13").Value = "Cooper (s) Hardtop & Cabrio"
$ws.Range("D14").Value = "Countryman"
$ws.Range("D15").Value = "Coupe"
$ws.Range("D16").Value = "John Cooper Works"

# --- Row 17: Year, Model, and Additional Information change ---
$ws.Range("B17").Value = 2013
$ws.Range("D17").Value = "Clubman"
$ws.Range("F17").Value = "Information not available"

# --- Update column D (Model) for rows 18-21 ---
$ws.Range("D18").Value = "Cooper (s) Hardtop & Cabrio"
$ws.Range("D19").Value = "Countryman"
$ws.Range("D20").Value = "Coupe"
$ws.Range("D21").Value = "John Cooper Works"

# --- Remove the now-redundant trailing rows 23 and 22 ---
# (delete from the bottom up so row numbers don't shift mid-operation)
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
